# Atualizacao de bases das ligas, do dia: 24-02-2024 as 23:13
# Swap the betting-odds data rows that were re-ordered upstream, and
# refresh a handful of closing-odds values for still-unplayed fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 201 and row 202: columns B,F,G,H,J,K,M,N,O,P,Q,R,S,T,U,V,W,X,Z,AA,AB,AC
$ws.Range("B201").Value = 6810777
$ws.Range("B202").Value = 6811260
$ws.Range("F201").Value = "Yverdon Sport FC"
$ws.Range("F202").Value = "FC Zurich"
$ws.Range("G201").Value = "Stade LausanneOuchy"
$ws.Range("G202").Value = "Lucerne"
$ws.Range("H201").Value = 2
$ws.Range("H202").Value = 1
$ws.Range("J201").Value = "H"
$ws.Range("J202").Value = "D"
$ws.Range("K201").Value = 2.3
$ws.Range("K202").Value = 1.833
$ws.Range("M201").Value = 2.875
$ws.Range("M202").Value = 4.2
$ws.Range("N201").Value = 3
$ws.Range("N202").Value = 1.7
$ws.Range("O201").Value = 3.5
$ws.Range("O202").Value = 3.6
$ws.Range("P201").Value = 2.375
$ws.Range("P202").Value = 4.75
$ws.Range("Q201").Value = 0.25
$ws.Range("Q202").Value = -0.75
$ws.Range("R201").Value = 1.8
$ws.Range("R202").Value = 1.95
$ws.Range("S201").Value = 2.05
$ws.Range("S202").Value = 1.9
$ws.Range("T201").Value = 2.5
$ws.Range("T202").Value = 2.75
$ws.Range("U201").Value = 1.925
$ws.Range("U202").Value = 1.975
$ws.Range("V201").Value = 1.925
$ws.Range("V202").Value = 1.875
$ws.Range("W201").Value = 2
$ws.Range("W202").Value = -1
$ws.Range("X201").Value = -1
$ws.Range("X202").Value = 2.6
$ws.Range("Z201").Value = 0.8
$ws.Range("Z202").Value = -1
$ws.Range("AA201").Value = -1
$ws.Range("AA202").Value = 0.8999999999999999
$ws.Range("AB201").Value = 0.925
$ws.Range("AB202").Value = -1
$ws.Range("AC201").Value = -1
$ws.Range("AC202").Value = 0.875

# Swap row 214 and row 215: columns B,F,G,H,I,J,K,M,O,P,R,S,U,V,W,X,Z,AA,AB,AC
$ws.Range("B214").Value = 6810779
$ws.Range("B215").Value = 6810310
$ws.Range("F214").Value = "Yverdon Sport FC"
$ws.Range("F215").Value = "Winterthur"
$ws.Range("G214").Value = "Lucerne"
$ws.Range("G215").Value = "Servette"
$ws.Range("H214").Value = 2
$ws.Range("H215").Value = 3
$ws.Range("I214").Value = 1
$ws.Range("I215").Value = 3
$ws.Range("J214").Value = "H"
$ws.Range("J215").Value = "D"
$ws.Range("K214").Value = 2.875
$ws.Range("K215").Value = 3.5
$ws.Range("M214").Value = 2.3
$ws.Range("M215").Value = 2
$ws.Range("O214").Value = 3.75
$ws.Range("O215").Value = 3.6
$ws.Range("P214").Value = 2.1
$ws.Range("P215").Value = 2.15
$ws.Range("R214").Value = 2
$ws.Range("R215").Value = 1.95
$ws.Range("S214").Value = 1.85
$ws.Range("S215").Value = 1.9
$ws.Range("U214").Value = 1.8
$ws.Range("U215").Value = 1.875
$ws.Range("V214").Value = 2.05
$ws.Range("V215").Value = 1.975
$ws.Range("W214").Value = 2.2
$ws.Range("W215").Value = -1
$ws.Range("X214").Value = -1
$ws.Range("X215").Value = 2.6
$ws.Range("Z214").Value = 1
$ws.Range("Z215").Value = 0.475
$ws.Range("AA214").Value = -1
$ws.Range("AA215").Value = -0.5
$ws.Range("AB214").Value = 0.4
$ws.Range("AB215").Value = 0.875
$ws.Range("AC214").Value = -0.5
$ws.Range("AC215").Value = -1

# Swap row 220 and row 221: columns B,F,G,K,L,M,N,O,P,Q,R,S,T,U,V,W,Z,AB,AC
$ws.Range("B220").Value = 6811273
$ws.Range("B221").Value = 6811272
$ws.Range("F220").Value = "Lucerne"
$ws.Range("F221").Value = "Grasshoppers"
$ws.Range("G220").Value = "Lausanne Sports"
$ws.Range("G221").Value = "FC Zurich"
$ws.Range("K220").Value = 1.727
$ws.Range("K221").Value = 3.75
$ws.Range("L220").Value = 3.5
$ws.Range("L221").Value = 3.6
$ws.Range("M220").Value = 5
$ws.Range("M221").Value = 1.909
$ws.Range("N220").Value = 2.05
$ws.Range("N221").Value = 3.5
$ws.Range("O220").Value = 3.6
$ws.Range("O221").Value = 3.8
$ws.Range("P220").Value = 3.5
$ws.Range("P221").Value = 1.95
$ws.Range("Q220").Value = -0.5
$ws.Range("Q221").Value = 0.5
$ws.Range("R220").Value = 2.025
$ws.Range("R221").Value = 1.925
$ws.Range("S220").Value = 1.825
$ws.Range("S221").Value = 1.925
$ws.Range("T220").Value = 2.75
$ws.Range("T221").Value = 2.5
$ws.Range("U220").Value = 1.875
$ws.Range("U221").Value = 1.925
$ws.Range("V220").Value = 1.975
$ws.Range("V221").Value = 1.925
$ws.Range("W220").Value = 1.05
$ws.Range("W221").Value = 2.5
$ws.Range("Z220").Value = 1.025
$ws.Range("Z221").Value = 0.925
$ws.Range("AB220").Value = 0.4375
$ws.Range("AB221").Value = 0.925
$ws.Range("AC220").Value = -0.5
$ws.Range("AC221").Value = -1

# Swap row 238 and row 239: columns B,F,G,H,I,K,L,M,N,O,P,Q,R,S,T,U,V,W,Z,AB,AC
$ws.Range("B238").Value = 7616924
$ws.Range("B239").Value = 7616836
$ws.Range("F238").Value = "Basel"
$ws.Range("F239").Value = "Yverdon Sport FC"
$ws.Range("G238").Value = "St Gallen"
$ws.Range("G239").Value = "Servette"
$ws.Range("H238").Value = 1
$ws.Range("H239").Value = 2
$ws.Range("I238").Value = 0
$ws.Range("I239").Value = 1
$ws.Range("K238").Value = 2.875
$ws.Range("K239").Value = 4.5
$ws.Range("L238").Value = 3.75
$ws.Range("L239").Value = 4.333
$ws.Range("M238").Value = 2.25
$ws.Range("M239").Value = 1.615
$ws.Range("N238").Value = 2.3
$ws.Range("N239").Value = 4
$ws.Range("O238").Value = 3.5
$ws.Range("O239").Value = 3.75
$ws.Range("P238").Value = 3.1
$ws.Range("P239").Value = 1.85
$ws.Range("Q238").Value = -0.25
$ws.Range("Q239").Value = 0.5
$ws.Range("R238").Value = 2
$ws.Range("R239").Value = 1.95
$ws.Range("S238").Value = 1.85
$ws.Range("S239").Value = 1.9
$ws.Range("T238").Value = 2.75
$ws.Range("T239").Value = 2.5
$ws.Range("U238").Value = 2.025
$ws.Range("U239").Value = 1.85
$ws.Range("V238").Value = 1.825
$ws.Range("V239").Value = 2
$ws.Range("W238").Value = 1.3
$ws.Range("W239").Value = 3
$ws.Range("Z238").Value = 1
$ws.Range("Z239").Value = 0.95
$ws.Range("AB238").Value = -1
$ws.Range("AB239").Value = 0.8500000000000001
$ws.Range("AC238").Value = 0.825
$ws.Range("AC239").Value = -1


# Row 247 (id 247) - refreshed closing odds
$ws.Range("O249").Value = 3.2
$ws.Range("P249").Value = 2.7
$ws.Range("T249").Value = 2.25
$ws.Range("U249").Value = 1.8
$ws.Range("V249").Value = 2.05

# Row 249 (id 249) - refreshed closing odds
$ws.Range("R251").Value = 1.875
$ws.Range("S251").Value = 1.975
